$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no auto numeric conversion) for D-column cells
# whose new price values would otherwise be parsed as numbers by Excel.
$textCells = @("D4","D5","D8","D9","D10","D11","D12","D13","D14","D15","D16","D18","D19","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D36","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = "29.403.10"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.849.61"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "240.49"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.07641"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "0.2907"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").Value = "24.79"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "0.07744"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "5.039"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "0.6799"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "0.00001070"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "83.32"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "6.172"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "29.432.15"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "227.95"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "12.34"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "7.470"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "157.98"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "0.1383"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "8.426"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "17.70"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "1.382"
$ws.Range("E27").Value = "  +6.11%  "
$ws.Range("D28").Value = "1.464"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "0.05597"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "4.131"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "4.064"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").Value = "1.841"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "1.164"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "0.6963"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "0.01801"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "1.229.73"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "2.722"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "6.411"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "0.9083"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "101.67"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "65.96"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "7.199"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "0.00000000120"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "0.4016"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "8.995"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "1.683"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "0.1146"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").Value = "0.05704"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "0.4633"
$ws.Range("E51").Value = "  +0.17%  "
